# Updates cryptos list values per the scheduled GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '70.172.98'
$ws.Range("E2").Value = '  +1.19%  '

# Row 3
$ws.Range("D3").Value = '3.506.55'
$ws.Range("E3").Value = '  -0.43%  '

# Row 4
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.32'
$ws.Range("E5").Value = '  -0.69%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.61'
$ws.Range("E6").Value = '  +2.00%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.610'
$ws.Range("E7").Value = '  -0.57%  '

# Row 8
$ws.Range("D8").Value = '3.500.48'
$ws.Range("E8").Value = '  -0.49%  '

# Row 9
$ws.Range("E9").Value = '  +0.02%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.193'
$ws.Range("E10").Value = '  -1.39%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.23'
$ws.Range("E11").Value = '  +8.60%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.582'
$ws.Range("E12").Value = '  +0.21%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.10'
$ws.Range("E13").Value = '  -2.65%  '

# Row 14
$ws.Range("E14").Value = '  -1.39%  '

# Row 15
$ws.Range("D15").Value = '4.076.80'
$ws.Range("E15").Value = '  -0.60%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.27'
$ws.Range("E16").Value = '  -1.15%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '609.79'
$ws.Range("E17").Value = '  -0.90%  '

# Row 18
$ws.Range("D18").Value = '3.516.56'
$ws.Range("E18").Value = '  -0.40%  '

# Row 19
$ws.Range("D19").Value = '70.318.77'
$ws.Range("E19").Value = '  +1.18%  '

# Row 20
$ws.Range("E20").Value = '  +0.58%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.28'
$ws.Range("E21").Value = '  +0.09%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.872'
$ws.Range("E22").Value = '  -1.18%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.03'
$ws.Range("E23").Value = '  -19.68%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.61'
$ws.Range("E24").Value = '  -1.41%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.42'
$ws.Range("E25").Value = '  +0.90%  '

# Row 26
$ws.Range("E26").Value = '  -4.35%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.12%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.55'
$ws.Range("E28").Value = '  -3.19%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.98'
$ws.Range("E29").Value = '  +2.27%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.94'
$ws.Range("E30").Value = '  -3.64%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.09'
$ws.Range("E31").Value = '  -4.78%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.97'
$ws.Range("E32").Value = '  -4.84%  '

# Row 33
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.87'
$ws.Range("E33").Value = '  -1.04%  '

# Row 34
$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '634.33'
$ws.Range("E34").Value = '  +11.09%  '

# Row 35
$ws.Range("E35").Value = '  -5.13%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.57'
$ws.Range("E36").Value = '  +0.49%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0993'
$ws.Range("E37").Value = '  -2.41%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.69'
$ws.Range("E38").Value = '  -0.98%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0471'
$ws.Range("E39").Value = '  +5.86%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '56.64'
$ws.Range("E40").Value = '  -0.60%  '

# Row 41
$ws.Range("E41").Value = '  +0.01%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.141'
$ws.Range("E42").Value = '  +1.21%  '

# Row 43
$ws.Range("D43").Value = '0.0₃0737'
$ws.Range("E43").Value = '  +4.56%  '

# Row 44
$ws.Range("D44").Value = '3.351.86'
$ws.Range("E44").Value = '  -1.02%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.309'
$ws.Range("E45").Value = '  -5.55%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.89'
$ws.Range("E46").Value = '  -0.18%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '32.07'
$ws.Range("E47").Value = '  -2.99%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.54'
$ws.Range("E48").Value = '  -2.23%  '

# Row 49
$ws.Range("E49").Value = '  +0.03%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.85'
$ws.Range("E50").Value = '  -0.20%  '

# Row 51
$ws.Range("E51").Value = '  -0.01%  '
